$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E10").Value = 527
$ws.Range("E11").Value = 341
$ws.Range("E12").Value = 518
$ws.Range("F12").Value = 279
$ws.Range("H12").Value = 364
$ws.Range("E13").Value = 130
$ws.Range("E15").Value = 163
$ws.Range("F15").Value = 72
$ws.Range("H15").Value = 121
$ws.Range("E16").Value = 201
$ws.Range("F16").Value = 101
$ws.Range("H16").Value = 149
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 54
$ws.Range("H17").Value = 78
$ws.Range("E24").Value = 209
$ws.Range("E27").Value = 321
$ws.Range("F27").Value = 161
$ws.Range("H27").Value = 242
$ws.Range("E29").Value = 166
$ws.Range("F29").Value = 93
$ws.Range("H29").Value = 134
$ws.Range("E30").Value = 203
$ws.Range("E31").Value = 73
$ws.Range("E32").Value = 178
$ws.Range("F32").Value = 105
$ws.Range("H32").Value = 143
$ws.Range("E34").Value = 213
$ws.Range("E37").Value = 156
$ws.Range("F45").Value = 70
$ws.Range("H45").Value = 109
$ws.Range("E46").Value = 314
$ws.Range("E47").Value = 447
$ws.Range("E48").Value = 204
$ws.Range("F49").Value = 121
$ws.Range("H49").Value = 208
